# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback DateTime" (column H) timestamps on the
# zh-cn and de-de worksheets (row 2) to reflect the new handback run.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-19 08:34:38"
$wsZh.Range("H2").Value = "2016-03-19 08:35:20"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-19 08:34:45"
$wsDe.Range("H2").Value = "2016-03-19 08:35:34"
